# Correct ontologies and stories naming
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "HighSugars"
$ws.Range("F3").Value = "Cancer"

$ws.Range("D4").Value = "Obesity"
$ws.Range("F4").Value = "Cancer"

$ws.Range("D7").Value = "Sport AND UsesBodyPart VALUE UpperBodyPart"

$ws.Range("D8").Value = "Swimming"

$ws.Range("E9").Value = "causesNutrientState"

$ws.Range("D10").Value = "Beef, Lamb"
$ws.Range("E10").Value = "causesNutrientState"

$ws.Range("D11").Value = "Kimbap, Rabokki"
$ws.Range("F11").Value = "KoreanFood"

$ws.Range("D12").Value = "Carbonara"
$ws.Range("E12").Value = "contains"
$ws.Range("F12").Value = "Cream"

$ws.Range("D13").Value = "Nuts, Strawberry, Blueberries"

$ws.Range("D14").Value = "Nuts, OliveOil, CoconutOil"
$ws.Range("F14").Value = "Cancer, Alzheimer's"

$ws.Range("D15").Value = "Fish, WholeGrains, GreenLeafyVegetables, Olives, Nuts"

$ws.Range("D16").Value = "Alcohol, Obesity, NoSport, Tobacco"
$ws.Range("F16").Value = "Cancer"

$ws.Range("D17").Value = "LowVitamin"
$ws.Range("F17").Value = "HairLoss"
